# Auto-generated edit script: refresh market-price derived columns (H,I,J,K,L,M,N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 27682
$ws.Range("J57").Value = 27682
$ws.Range("L57").Value = 83046
$ws.Range("N57").Value = -84044
$ws.Range("H70").Value = 1499.8572
$ws.Range("J70").Value = 1583.1666
$ws.Range("L70").Value = 4749.4998
$ws.Range("N70").Value = -5289.4998
$ws.Range("H73").Value = 1499.8572
$ws.Range("J73").Value = 1583.1666
$ws.Range("L73").Value = 4749.4998
$ws.Range("N73").Value = -6621.4998
$ws.Range("H132").Value = 1634.1526
$ws.Range("I132").Value = 1428.3508
$ws.Range("K132").Value = 4285.0524
$ws.Range("M132").Value = -1755.0524
$ws.Range("H135").Value = 973.3333
$ws.Range("I135").Value = 603.0625
$ws.Range("J135").Value = 2666
$ws.Range("K135").Value = 5427.5625
$ws.Range("L135").Value = 23994
$ws.Range("M135").Value = -2892.5625
$ws.Range("N135").Value = -29064
$ws.Range("H138").Value = 1847.4231
$ws.Range("J138").Value = 3819.2
$ws.Range("L138").Value = 11457.6
$ws.Range("N138").Value = -21737.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 518989.84
$ws.Range("I32").Value = 587196.2
$ws.Range("J32").Value = 18809.777
$ws.Range("K32").Value = 587196.2
$ws.Range("L32").Value = 18809.777
$ws.Range("M32").Value = -586909.2
$ws.Range("N32").Value = -19383.777
$ws.Range("H61").Value = 2370.7903
$ws.Range("I61").Value = 1966.8096
$ws.Range("J61").Value = 3219.15
$ws.Range("K61").Value = 1966.8096
$ws.Range("L61").Value = 3219.15
$ws.Range("M61").Value = -1754.8096
$ws.Range("N61").Value = -3643.15
$ws.Range("H74").Value = 2387.1943
$ws.Range("I74").Value = 2106.2307
$ws.Range("J74").Value = 3117.7
$ws.Range("K74").Value = 2106.2307
$ws.Range("L74").Value = 3117.7
$ws.Range("M74").Value = -1232.2307
$ws.Range("N74").Value = -4865.7
$ws.Range("H77").Value = 2387.1943
$ws.Range("I77").Value = 2106.2307
$ws.Range("J77").Value = 3117.7
$ws.Range("K77").Value = 10531.1535
$ws.Range("L77").Value = 15588.5
$ws.Range("M77").Value = -6163.1535
$ws.Range("N77").Value = -24324.5
$ws.Range("H136").Value = 2370.7903
$ws.Range("I136").Value = 1966.8096
$ws.Range("J136").Value = 3219.15
$ws.Range("K136").Value = 5900.4288
$ws.Range("L136").Value = 9657.450000000001
$ws.Range("M136").Value = -3350.4288
$ws.Range("N136").Value = -14757.45

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3308.6667
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 3308.6667
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H139").Value = 49890
$ws.Range("J139").Value = 49890
$ws.Range("L139").Value = 49890
$ws.Range("N139").Value = -60170

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 57476
$ws.Range("J52").Value = 57476
$ws.Range("L52").Value = 57476
$ws.Range("N52").Value = -58064
$ws.Range("H62").Value = 2540
$ws.Range("I62").Value = 2482.5
$ws.Range("K62").Value = 2482.5
$ws.Range("M62").Value = -1858.5
$ws.Range("H65").Value = 2540
$ws.Range("I65").Value = 2482.5
$ws.Range("K65").Value = 12412.5
$ws.Range("M65").Value = -9292.5
$ws.Range("H117").Value = 34575
$ws.Range("I117").Value = 19800
$ws.Range("J117").Value = 49350
$ws.Range("K117").Value = 19800
$ws.Range("L117").Value = 49350
$ws.Range("M117").Value = -15211
$ws.Range("N117").Value = -58528
$ws.Range("H125").Value = 38331.668
$ws.Range("J125").Value = 38331.668
$ws.Range("L125").Value = 38331.668
$ws.Range("N125").Value = -43251.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 1924.1666
$ws.Range("J31").Value = 1826.3636
$ws.Range("L31").Value = 5479.0908
$ws.Range("N31").Value = -6055.0908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 1000
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 1000
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -1504
$ws.Range("H36").Value = 2366.75
$ws.Range("I36").Value = 2233.5
$ws.Range("J36").Value = 2500
$ws.Range("K36").Value = 2233.5
$ws.Range("L36").Value = 2500
$ws.Range("M36").Value = -1748.5
$ws.Range("N36").Value = -3470
$ws.Range("H70").Value = 7902.5884
$ws.Range("I70").Value = 8572.27
$ws.Range("J70").Value = 5726.125
$ws.Range("K70").Value = 8572.27
$ws.Range("L70").Value = 5726.125
$ws.Range("M70").Value = -8302.27
$ws.Range("N70").Value = -6266.125
$ws.Range("H73").Value = 7902.5884
$ws.Range("I73").Value = 8572.27
$ws.Range("J73").Value = 5726.125
$ws.Range("K73").Value = 8572.27
$ws.Range("L73").Value = 5726.125
$ws.Range("M73").Value = -7636.27
$ws.Range("N73").Value = -7598.125
$ws.Range("H132").Value = 3267.0732
$ws.Range("I132").Value = 2942.75
$ws.Range("J132").Value = 5602.2
$ws.Range("K132").Value = 8828.25
$ws.Range("L132").Value = 16806.6
$ws.Range("M132").Value = -6298.25
$ws.Range("N132").Value = -21866.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 9525873
$ws.Range("I16").Value = 2151.2307
$ws.Range("J16").Value = 71430060
$ws.Range("K16").Value = 2151.2307
$ws.Range("L16").Value = 71430060
$ws.Range("M16").Value = -1981.2307
$ws.Range("N16").Value = -71430400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 225873.33
$ws.Range("J110").Value = 225873.33
$ws.Range("L110").Value = 225873.33
$ws.Range("N110").Value = -234053.33
$ws.Range("H111").Value = 17296.666
$ws.Range("J111").Value = 17296.666
$ws.Range("L111").Value = 17296.666
$ws.Range("N111").Value = -25476.666
$ws.Range("H117").Value = 24704.5
$ws.Range("J117").Value = 24704.5
$ws.Range("L117").Value = 24704.5
$ws.Range("N117").Value = -33882.5
$ws.Range("H124").Value = 36976.332
$ws.Range("J124").Value = 36976.332
$ws.Range("L124").Value = 36976.332
$ws.Range("N124").Value = -46796.332
$ws.Range("H127").Value = 54607.25
$ws.Range("I127").Value = 25000
$ws.Range("K127").Value = 25000
$ws.Range("M127").Value = -20040
